$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Articles")
$ws.Activate()

# ---- Row 16 ----
$ws.Range("A16").Value = 6
$ws.Range("B16").Formula = "=VLOOKUP(A16,Sites!A:C,3,0)"

$ws.Range("C15").Copy($ws.Range("C16"))
$ws.Range("C16").Value = "https://build-de.blogspot.com/2021/01/yuumi-build-playstyle-was-halten-die.html"
$ws.Hyperlinks.Add($ws.Range("D16"), "https://lolesports.com/article/%E2%80%9Cit-can-get-really-brutal%E2%80%9D-%E2%80%93-what-do-the-pros-think-of-yuumi/bltbc4c66cd3645937a") | Out-Null

$ws.Range("D15").Copy($ws.Range("D16"))
$ws.Range("D16").Value = "https://lolesports.com/article/%E2%80%9Cit-can-get-really-brutal%E2%80%9D-%E2%80%93-what-do-the-pros-think-of-yuumi/bltbc4c66cd3645937a"

$ws.Range("E15").Copy($ws.Range("E16"))
$ws.Range("E16").Value = 44203

# ---- Row 17 ----
$ws.Range("A17").Value = 6
$ws.Range("B17").Formula = "=VLOOKUP(A17,Sites!A:C,3,0)"

$ws.Range("D15").Copy($ws.Range("D17"))
$ws.Range("D17").Value = "https://lolesports.com/article/meta-reflection-into-the-jungle-with-sk-trick,-mad-shadow,-and-rge-inspired/blt7a7054fa7a9631af"

$ws.Range("E17").Value = "not posted yet"

# ---- Row 18 ----
$ws.Range("B18").Formula = "=VLOOKUP(A18,Sites!A:C,3,0)"

$ws.Range("D15").Copy($ws.Range("D18"))
$ws.Range("D18").Value = "https://lolesports.com/article/mad-kaiser-and-rge-vander-break-down-this-season's-support-meta/blte7822744fe3306f2"

$ws.Range("E18").Value = "not posted yet"

# ---- Row 19 ----
$ws.Range("B19").Formula = "=VLOOKUP(A19,Sites!A:C,3,0)"

$ws.Range("D15").Copy($ws.Range("D19"))
$ws.Range("D19").Value = "https://lolesports.com/article/%E2%80%9Ci-hope-the-meta-will-change-every-season%E2%80%9D-%E2%80%93-analysing-the-summer-season-mid-lane-with-mad-humanoid/blted96f359ec766296"

$ws.Range("E19").Value = "not posted yet"

# ---- Row 20 ----
$ws.Range("B20").Formula = "=VLOOKUP(A20,Sites!A:C,3,0)"

$ws.Range("D15").Copy($ws.Range("D20"))
$ws.Range("D20").Value = "https://lolesports.com/article/%E2%80%9Ci-liked-the-meta-more-than-in-spring%E2%80%9D-%E2%80%93-reflecting-on-summer-season%E2%80%99s-top-lane-with-mad-lions%E2%80%99-orome/blte7d733406124c06c"

$ws.Range("E20").Value = "not posted yet"

# ---- Hyperlinks (added in the same order as source rIds 17..22) ----
$ws.Hyperlinks.Add($ws.Range("D17"), "https://lolesports.com/article/meta-reflection-into-the-jungle-with-sk-trick,-mad-shadow,-and-rge-inspired/blt7a7054fa7a9631af") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D18"), "https://lolesports.com/article/mad-kaiser-and-rge-vander-break-down-this-season's-support-meta/blte7822744fe3306f2") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D19"), "https://lolesports.com/article/%E2%80%9Ci-hope-the-meta-will-change-every-season%E2%80%9D-%E2%80%93-analysing-the-summer-season-mid-lane-with-mad-humanoid/blted96f359ec766296") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D20"), "https://lolesports.com/article/%E2%80%9Ci-liked-the-meta-more-than-in-spring%E2%80%9D-%E2%80%93-reflecting-on-summer-season%E2%80%99s-top-lane-with-mad-lions%E2%80%99-orome/blte7d733406124c06c") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C16"), "https://build-de.blogspot.com/2021/01/yuumi-build-playstyle-was-halten-die.html") | Out-Null

# Restore the "Link" style on hyperlink cells (Hyperlinks.Add overwrites cell style)
$ws.Range("D15").Copy($ws.Range("D16"))
$ws.Range("D16").Value = "https://lolesports.com/article/%E2%80%9Cit-can-get-really-brutal%E2%80%9D-%E2%80%93-what-do-the-pros-think-of-yuumi/bltbc4c66cd3645937a"
$ws.Range("D15").Copy($ws.Range("D17"))
$ws.Range("D17").Value = "https://lolesports.com/article/meta-reflection-into-the-jungle-with-sk-trick,-mad-shadow,-and-rge-inspired/blt7a7054fa7a9631af"
$ws.Range("D15").Copy($ws.Range("D18"))
$ws.Range("D18").Value = "https://lolesports.com/article/mad-kaiser-and-rge-vander-break-down-this-season's-support-meta/blte7822744fe3306f2"
$ws.Range("D15").Copy($ws.Range("D19"))
$ws.Range("D19").Value = "https://lolesports.com/article/%E2%80%9Ci-hope-the-meta-will-change-every-season%E2%80%9D-%E2%80%93-analysing-the-summer-season-mid-lane-with-mad-humanoid/blted96f359ec766296"
$ws.Range("D15").Copy($ws.Range("D20"))
$ws.Range("D20").Value = "https://lolesports.com/article/%E2%80%9Ci-liked-the-meta-more-than-in-spring%E2%80%9D-%E2%80%93-reflecting-on-summer-season%E2%80%99s-top-lane-with-mad-lions%E2%80%99-orome/blte7d733406124c06c"
$ws.Range("C15").Copy($ws.Range("C16"))
$ws.Range("C16").Value = "https://build-de.blogspot.com/2021/01/yuumi-build-playstyle-was-halten-die.html"

$excel.Calculate()

$ws.Range("H18").Select()

Write-Host "done"
